$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-07 Friday" "2025-11-08 Saturday"

Replace-Text "918÷2=" "433÷9="
Replace-Text "990÷3=" "713÷5="
Replace-Text "922÷3=" "438÷3="
Replace-Text "835÷2=" "448÷7="
Replace-Text "170÷7=" "672÷6="
Replace-Text "980÷7=" "434÷2="
Replace-Text "561÷5=" "602÷4="
Replace-Text "402÷6=" "384÷4="
Replace-Text "391÷4=" "914÷7="
Replace-Text "501÷9=" "494÷4="
Replace-Text "584÷6=" "110÷6="
Replace-Text "319÷5=" "327÷6="
Replace-Text "450÷2=" "681÷8="
Replace-Text "823÷9=" "234÷6="
Replace-Text "701÷2=" "355÷9="
Replace-Text "134÷5=" "179÷9="
Replace-Text "631÷8=" "419÷3="
Replace-Text "111÷8=" "939÷3="
Replace-Text "962÷4=" "950÷7="
Replace-Text "151÷7=" "810÷2="
Replace-Text "710÷4=" "525÷2="
Replace-Text "492÷2=" "906÷7="
Replace-Text "550÷7=" "185÷5="
Replace-Text "326÷7=" "984÷9="
Replace-Text "790÷9=" "405÷2="
